$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 88
$ws.Range("H88").Value = 1587.7
$ws.Range("I88").Value = 1889.7142
$ws.Range("J88").Value = 883
$ws.Range("K88").Value = 1889.7142
$ws.Range("L88").Value = 883
$ws.Range("M88").Value = -1483.7142
$ws.Range("N88").Value = -1695
# Row 91
$ws.Range("H91").Value = 1587.7
$ws.Range("I91").Value = 1889.7142
$ws.Range("J91").Value = 883
$ws.Range("K91").Value = 1889.7142
$ws.Range("L91").Value = 883
$ws.Range("M91").Value = -485.7141999999999
$ws.Range("N91").Value = -3691
# Row 115
$ws.Range("H115").Value = 718.4286
$ws.Range("I115").Value = 718.4286
$ws.Range("K115").Value = 2155.2858
$ws.Range("M115").Value = -588.2857999999997
# Row 118
$ws.Range("H118").Value = 1553.8182
$ws.Range("I118").Value = 1632.6666
$ws.Range("J118").Value = 1199
$ws.Range("K118").Value = 4897.9998
$ws.Range("L118").Value = 3597
$ws.Range("M118").Value = -3240.9998
$ws.Range("N118").Value = -6911
# Row 137
$ws.Range("H137").Value = 2321.3142
$ws.Range("I137").Value = 2470.1538
$ws.Range("J137").Value = 1891.3334
$ws.Range("K137").Value = 7410.4614
$ws.Range("L137").Value = 5674.0002
$ws.Range("M137").Value = -4860.4614
$ws.Range("N137").Value = -10774.0002
# Row 138
$ws.Range("H138").Value = 8776735
$ws.Range("I138").Value = 1240.1578
$ws.Range("K138").Value = 3720.4734
$ws.Range("M138").Value = 1419.5266

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 539.08
$ws.Range("I2").Value = 411.47827
$ws.Range("J2").Value = 2006.5
$ws.Range("K2").Value = 411.47827
$ws.Range("L2").Value = 2006.5
$ws.Range("M2").Value = -298.47827
$ws.Range("N2").Value = -2232.5
# Row 32
$ws.Range("H32").Value = 10875350
$ws.Range("I32").Value = 14288064
$ws.Range("K32").Value = 14288064
$ws.Range("M32").Value = -14287777
# Row 74
$ws.Range("H74").Value = 62571760
$ws.Range("I74").Value = 77010296
$ws.Range("K74").Value = 77010296
$ws.Range("M74").Value = -77009422
# Row 77
$ws.Range("H77").Value = 62571760
$ws.Range("I77").Value = 77010296
$ws.Range("K77").Value = 385051480
$ws.Range("M77").Value = -385047112
# Row 116
$ws.Range("H116").Value = 539.08
$ws.Range("I116").Value = 411.47827
$ws.Range("J116").Value = 2006.5
$ws.Range("K116").Value = 411.47827
$ws.Range("L116").Value = 2006.5
$ws.Range("M116").Value = 1882.52173
$ws.Range("N116").Value = -6594.5
# Row 122
$ws.Range("H122").Value = 3845.3
$ws.Range("I122").Value = 3013.2856
$ws.Range("K122").Value = 9039.856800000001
$ws.Range("M122").Value = -6589.856800000001
# Row 141
$ws.Range("H141").Value = 47500
$ws.Range("J141").Value = 47500
$ws.Range("L141").Value = 47500
$ws.Range("N141").Value = -57860

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 539.08
$ws.Range("I3").Value = 411.47827
$ws.Range("J3").Value = 2006.5
$ws.Range("K3").Value = 411.47827
$ws.Range("L3").Value = 2006.5
$ws.Range("M3").Value = -297.47827
$ws.Range("N3").Value = -2234.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 23260292
$ws.Range("I31").Value = 3643.5667
$ws.Range("K31").Value = 3643.5667
$ws.Range("M31").Value = -3348.5667
# Row 34
$ws.Range("H34").Value = 23260292
$ws.Range("I34").Value = 3643.5667
$ws.Range("K34").Value = 3643.5667
$ws.Range("M34").Value = -3441.5667
# Row 58
$ws.Range("H58").Value = 2202.0435
$ws.Range("I58").Value = 2089.9412
$ws.Range("K58").Value = 2089.9412
$ws.Range("M58").Value = -1886.9412
# Row 99
$ws.Range("H99").Value = 16830.842
$ws.Range("I99").Value = 26595.334
$ws.Range("K99").Value = 26595.334
$ws.Range("M99").Value = -25097.334
# Row 126
$ws.Range("H126").Value = 16830.842
$ws.Range("I126").Value = 26595.334
$ws.Range("K126").Value = 79786.00199999999
$ws.Range("M126").Value = -77316.00199999999
# Row 134
$ws.Range("H134").Value = 1330.96
$ws.Range("I134").Value = 1185.826
$ws.Range("K134").Value = 3557.478
$ws.Range("M134").Value = -1022.478
# Row 136
$ws.Range("H136").Value = 2202.0435
$ws.Range("I136").Value = 2089.9412
$ws.Range("K136").Value = 6269.823600000001
$ws.Range("M136").Value = -3719.823600000001
# Row 141
$ws.Range("H141").Value = 402938.34
$ws.Range("J141").Value = 512218.34
$ws.Range("L141").Value = 512218.34
$ws.Range("N141").Value = -522578.34

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 40707230
$ws.Range("J4").Value = 16795812
$ws.Range("L4").Value = 50387436
$ws.Range("N4").Value = -50387660
# Row 128
$ws.Range("H128").Value = 115994.75
$ws.Range("I128").Value = 115994.75
$ws.Range("K128").Value = 347984.25
$ws.Range("M128").Value = -343004.25
# Row 131
$ws.Range("H131").Value = 32877.61
$ws.Range("J131").Value = 5214.5386
$ws.Range("L131").Value = 15643.6158
$ws.Range("N131").Value = -25723.6158
# Row 132
$ws.Range("H132").Value = 1757039.8
$ws.Range("I132").Value = 2185.476
$ws.Range("J132").Value = 3924801
$ws.Range("K132").Value = 19669.284
$ws.Range("L132").Value = 35323209
$ws.Range("M132").Value = -17139.284
$ws.Range("N132").Value = -35328269
# Row 134
$ws.Range("H134").Value = 2487.111
$ws.Range("I134").Value = 1456.9412
$ws.Range("K134").Value = 4370.8236
$ws.Range("M134").Value = 699.1764000000003

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 4629.25
$ws.Range("I43").Value = 1008.5
$ws.Range("J43").Value = 8250
$ws.Range("K43").Value = 1008.5
$ws.Range("L43").Value = 8250
$ws.Range("M43").Value = -857.5
$ws.Range("N43").Value = -8552

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4041.2856
$ws.Range("I7").Value = 3563.3125
$ws.Range("J7").Value = 4678.5835
$ws.Range("K7").Value = 3563.3125
$ws.Range("L7").Value = 4678.5835
$ws.Range("M7").Value = -3451.3125
$ws.Range("N7").Value = -4902.5835
# Row 122
$ws.Range("H122").Value = 4643.2856
$ws.Range("I122").Value = 4000.6667
$ws.Range("K122").Value = 12002.0001
$ws.Range("M122").Value = -9552.000100000001
# Row 126
$ws.Range("H126").Value = 4041.2856
$ws.Range("I126").Value = 3563.3125
$ws.Range("J126").Value = 4678.5835
$ws.Range("K126").Value = 10689.9375
$ws.Range("L126").Value = 14035.7505
$ws.Range("M126").Value = -8219.9375
$ws.Range("N126").Value = -18975.7505
# Row 131
$ws.Range("H131").Value = 89778
$ws.Range("J131").Value = 89778
$ws.Range("L131").Value = 89778
$ws.Range("N131").Value = -99858
